$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "doFilterSearchInContacts"

$newSheet.Range("A1").Value = "value"
$newSheet.Range("A2").Value = "Dalal"
$newSheet.Range("A3").Value = "Abdullah"
$newSheet.Range("A4").Value = "Auooz"

$newSheet.Range("A5").Select() | Out-Null

# Update selection on the doSignIn sheet
$ws1 = $wb.Worksheets.Item("doSignIn")
$ws1.Activate()
$ws1.Range("A2").Select() | Out-Null
